$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Home")

# New values for existing rows 3-7 (column A stays the same; B/C/D filled in)
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 100
$ws.Range("D3").Value = 1800

$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 110
$ws.Range("D4").Value = 2200

$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 90
$ws.Range("D5").Value = 2000

$ws.Range("B6").HorizontalAlignment = -4108
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 100
$ws.Range("D6").Value = 1900

$ws.Range("B7").Value = 2
$ws.Range("C7").Value = 130
$ws.Range("D7").Value = 2000

# New rows 8-11
$ws.Range("A8").Value = 15000
$ws.Range("B8").Value = 3
$ws.Range("C8").Value = 120
$ws.Range("D8").Value = 2500

$ws.Range("A9").Value = 21000
$ws.Range("B9").Value = 5
$ws.Range("C9").Value = 120
$ws.Range("D9").Value = 2300

$ws.Range("A10").Value = 18000
$ws.Range("B10").Value = 4
$ws.Range("C10").Value = 110
$ws.Range("D10").Value = 2200

$ws.Range("A11").Value = 15000
$ws.Range("B11").Value = 3
$ws.Range("C11").Value = 120
$ws.Range("D11").Value = 3000

$ws.Range("D11").Select()
